$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.481.97'
$ws.Range('E2').Value = '  +1.16%  '

$ws.Range('D3').Value = '1.904.21'
$ws.Range('E3').Value = '  +2.44%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.73'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.649'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.22%  '

$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.94'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.343'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0707'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0999'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.84%  '

$ws.Range('D12').Value = '2.180.25'
$ws.Range('E12').Value = '  +2.42%  '

$ws.Range('E13').Value = '  +8.88%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.698'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.94%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.907.07'
$ws.Range('E15').Value = '  +2.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.83'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.69%  '

$ws.Range('D17').Value = '35.472.50'
$ws.Range('E17').Value = '  +1.17%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '72.06'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.47%  '

$ws.Range('D19').Value = '0.0₃0831'
$ws.Range('E19').Value = '  +4.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '243.20'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.65'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.72%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.84'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.73%  '

$ws.Range('E23').Value = '  +0.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.26'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +16.40%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.60'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.53'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.62%  '

$ws.Range('E28').Value = '  +1.50%  '

$ws.Range('E29').Value = '  +2.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.961'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +22.96%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.14'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.74%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0571'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.18%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.20'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.94%  '

$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.01'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.14%  '

$ws.Range('E35').Value = '  +8.72%  '

$ws.Range('E36').Value = '  -0.06%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.33'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.93%  '

$ws.Range('E38').Value = '  +2.68%  '

$ws.Range('E39').Value = '  +16.90%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '91.83'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.08%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0204'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.71%  '

$ws.Range('E42').Value = '  +5.35%  '

$ws.Range('D43').Value = '1.347.43'
$ws.Range('E43').Value = '  -0.48%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '49.17'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +41.21%  '

$ws.Range('E45').Value = '  +1.96%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.63'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.02%  '

$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('E48').Value = '  +0.19%  '

$ws.Range('E49').Value = '  +3.74%  '

$ws.Range('D50').Value = '2.091.45'
$ws.Range('E50').Value = '  +2.46%  '

$ws.Range('E51').Value = '  +1.83%  '
